# "About page/impt. discoveries: content has gone live"
#
# 1) Row 20 (About Page changes / "Write copy for achievements"): mark DONE
#    (style like the other completed rows) and update hours spent; add a note.
# 2) Row 24 ("Add to About page"): mark DONE and set hours spent to 0.
# 3) Insert a brand-new task row right after it: "Share database with
#    faculty for edits" (1 hour), pushing the "Important Discoveries"
#    block and the totals row down by one.
# 4) Grow Table1 to cover the new row and fix the selection/top-left view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row first, while everything below is still at its
#     original address, so later edits can target final row numbers. ---
$ws.Rows.Item(25).Insert()

# New row 25: "Share database with faculty for edits"
$ws.Range("B25").Value = "Share database with faculty for edits"
$ws.Range("C25").Value = 1
$ws.Rows.Item(25).RowHeight = 29

# Give the new B25 the same "in progress" look B20 currently has (fill +
# wrap), by copying B20's format before B20 itself is changed below.
$ws.Range("B20").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 20: content went live -> mark done like rows 17-19/21-23, bump
#     time spent, and record a note about the work with Christoph. ---
$ws.Range("B17").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = "Worked with Christoph to develop db content"

# --- Row 24: content went live -> mark done, 0 hours spent so far. ---
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D24").Value = 0

# --- Grow the table (ListObject) so the new row is part of Table1. ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E38"))

# --- Sheet view: scroll/selection moved along with the edits. ---
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("B25").Select()
